# Updated SFI GM, VR and UF
# - AP_Value (col N, row 2): advance appointment slot to 10:50-10:54 AM
# - DT_Value (col AB, row 2): advance date value to 06/11/2023
# - VR_Value (col AR, row 2): newly populated voice-record file name
# - UF_Prefix (col AU, row 2): advance uploaded form file name to 06112023

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "Appointment Date : 06/11/2023, Time : [ 10:50 AM to 10:54 AM ]"

# AB2 holds a date-like label that must stay plain text (it was stored as a
# shared string, not a real date serial) -- force text entry with a leading
# apostrophe, then reapply the Normal style so no stray number format sticks.
$ws.Range("AB2").Value = "'06/11/2023"
$ws.Range("AB2").Style = "Normal"

$ws.Range("AR2").Value = "voice_record_06112023"
$ws.Range("AU2").Value = "formshow_06112023"
